$d = $word.ActiveDocument

$replacements = @(
    @("37×40=", "66×77="),
    @("87×13=", "64×78="),
    @("91×93=", "19×48="),
    @("27×58=", "57×43="),
    @("17×79=", "21×93="),
    @("89×75=", "71×76="),
    @("94×18=", "47×14="),
    @("97×53=", "44×69="),
    @("44×61=", "36×93="),
    @("53×49=", "25×60="),
    @("64×68=", "69×71="),
    @("31×92=", "69×40="),
    @("76×76=", "18×80="),
    @("85×36=", "29×61="),
    @("80×31=", "73×88="),
    @("29×73=", "99×32="),
    @("21×92=", "11×98="),
    @("55×95=", "56×47="),
    @("63×93=", "30×53="),
    @("81×54=", "20×43="),
    @("11×24=", "58×74="),
    @("48×80=", "39×79="),
    @("76×85=", "33×16="),
    @("20×13=", "86×14="),
    @("65×13=", "67×25=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
